$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 503.18182
$ws.Range("I4").Value = 310.2857
$ws.Range("J4").Value = 840.75
$ws.Range("K4").Value = 310.2857
$ws.Range("L4").Value = 840.75
$ws.Range("M4").Value = -196.2857
$ws.Range("N4").Value = -1068.75
$ws.Range("H69").Value = 9337.361000000001
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 9337.361000000001
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 28012.083
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -29760.083
$ws.Range("H72").Value = 9337.361000000001
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 9337.361000000001
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 84036.24900000001
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -92772.24900000001
$ws.Range("H98").Value = 2214.2144
$ws.Range("I98").Value = 1750
$ws.Range("K98").Value = 1750
$ws.Range("M98").Value = -252
$ws.Range("H107").Value = 1829.7949
$ws.Range("I107").Value = 1884.5428
$ws.Range("K107").Value = 1884.5428
$ws.Range("M107").Value = 35.45720000000006
$ws.Range("H111").Value = 2904
$ws.Range("I111").Value = 2243
$ws.Range("K111").Value = 6729
$ws.Range("M111").Value = -3662
$ws.Range("H113").Value = 34600
$ws.Range("I113").Value = 34600
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 34600
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -31346
$ws.Range("N113").ClearContents()
$ws.Range("H122").Value = 2214.2144
$ws.Range("I122").Value = 1750
$ws.Range("K122").Value = 5250
$ws.Range("M122").Value = -2800
$ws.Range("H126").Value = 95779.5
$ws.Range("J126").Value = 95779.5
$ws.Range("L126").Value = 95779.5
$ws.Range("N126").Value = -105659.5
$ws.Range("H127").Value = 928.8
$ws.Range("I127").Value = 780.9286
$ws.Range("K127").Value = 2342.7858
$ws.Range("M127").Value = 2617.2142
$ws.Range("H128").Value = 137390
$ws.Range("J128").Value = 137390
$ws.Range("L128").Value = 137390
$ws.Range("N128").Value = -147350
$ws.Range("H129").Value = 1849.875
$ws.Range("H132").Value = 2380.2437
$ws.Range("I132").Value = 2395.6494
$ws.Range("J132").Value = 1194
$ws.Range("K132").Value = 7186.948199999999
$ws.Range("L132").Value = 3582
$ws.Range("M132").Value = -4656.948199999999
$ws.Range("N132").Value = -8642
$ws.Range("H133").Value = 110000
$ws.Range("J133").Value = 110000
$ws.Range("L133").Value = 110000
$ws.Range("N133").Value = -120120
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H137").Value = 4558.143
$ws.Range("I137").Value = 3651.1667
$ws.Range("J137").Value = 10000
$ws.Range("K137").Value = 10953.5001
$ws.Range("L137").Value = 30000
$ws.Range("M137").Value = -8403.500100000001
$ws.Range("N137").Value = -35100
$ws.Range("H138").Value = 3569.2324
$ws.Range("I138").Value = 1783.7693
$ws.Range("J138").Value = 3839.128
$ws.Range("K138").Value = 5351.3079
$ws.Range("L138").Value = 11517.384
$ws.Range("M138").Value = -211.3078999999998
$ws.Range("N138").Value = -21797.384
$ws.Range("H141").Value = 551.96155
$ws.Range("I141").Value = 551.96155
$ws.Range("K141").Value = 1655.88465
$ws.Range("M141").Value = 3524.11535

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18048830
$ws.Range("I32").Value = 18718422
$ws.Range("J32").Value = 11910911
$ws.Range("K32").Value = 18718422
$ws.Range("L32").Value = 11910911
$ws.Range("M32").Value = -18718135
$ws.Range("N32").Value = -11911485
$ws.Range("H45").Value = 4176.1
$ws.Range("I45").Value = 3600.5334
$ws.Range("K45").Value = 3600.5334
$ws.Range("M45").Value = -3223.5334
$ws.Range("H61").Value = 4089.4443
$ws.Range("I61").Value = 3950.8333
$ws.Range("J61").Value = 4366.6665
$ws.Range("K61").Value = 3950.8333
$ws.Range("L61").Value = 4366.6665
$ws.Range("M61").Value = -3738.8333
$ws.Range("N61").Value = -4790.6665
$ws.Range("H132").Value = 3564.7144
$ws.Range("I132").Value = 3155.3333
$ws.Range("J132").Value = 4301.6
$ws.Range("K132").Value = 9465.999899999999
$ws.Range("L132").Value = 12904.8
$ws.Range("M132").Value = -6935.999899999999
$ws.Range("N132").Value = -17964.8
$ws.Range("H136").Value = 4089.4443
$ws.Range("I136").Value = 3950.8333
$ws.Range("J136").Value = 4366.6665
$ws.Range("K136").Value = 11852.4999
$ws.Range("L136").Value = 13099.9995
$ws.Range("M136").Value = -9302.499899999999
$ws.Range("N136").Value = -18199.9995

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("N58").ClearContents()
$ws.Range("H94").Value = 1537.4
$ws.Range("I94").Value = 950.3570999999999
$ws.Range("J94").Value = 2284.5454
$ws.Range("K94").Value = 950.3570999999999
$ws.Range("L94").Value = 2284.5454
$ws.Range("M94").Value = -499.3570999999999
$ws.Range("N94").Value = -3186.5454
$ws.Range("H97").Value = 56714
$ws.Range("I97").Value = 38285.332
$ws.Range("J97").Value = 112000
$ws.Range("K97").Value = 38285.332
$ws.Range("L97").Value = 112000
$ws.Range("M97").Value = -37294.332
$ws.Range("N97").Value = -113982
$ws.Range("H105").Value = 2838.8125
$ws.Range("J105").Value = 3599
$ws.Range("L105").Value = 3599
$ws.Range("N105").Value = -7093
$ws.Range("H132").Value = 118984.5
$ws.Range("J132").Value = 118984.5
$ws.Range("L132").Value = 118984.5
$ws.Range("N132").Value = -129104.5
$ws.Range("H134").Value = 2859997.5
$ws.Range("I134").Value = 3403260
$ws.Range("K134").Value = 10209780
$ws.Range("M134").Value = -10207245

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3375
$ws.Range("J16").Value = 3750
$ws.Range("L16").Value = 3750
$ws.Range("N16").Value = -4324
$ws.Range("H31").Value = 3740.4666
$ws.Range("I31").Value = 2685.5833
$ws.Range("J31").Value = 4443.722
$ws.Range("K31").Value = 2685.5833
$ws.Range("L31").Value = 4443.722
$ws.Range("M31").Value = -2390.5833
$ws.Range("N31").Value = -5033.722
$ws.Range("H34").Value = 3740.4666
$ws.Range("I34").Value = 2685.5833
$ws.Range("J34").Value = 4443.722
$ws.Range("K34").Value = 2685.5833
$ws.Range("L34").Value = 4443.722
$ws.Range("M34").Value = -2483.5833
$ws.Range("N34").Value = -4847.722
$ws.Range("H58").Value = 3119.2856
$ws.Range("I58").Value = 3015.3157
$ws.Range("J58").Value = 4107
$ws.Range("K58").Value = 3015.3157
$ws.Range("L58").Value = 4107
$ws.Range("M58").Value = -2812.3157
$ws.Range("N58").Value = -4513
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").ClearContents()
$ws.Range("H113").Value = 3375
$ws.Range("J113").Value = 3750
$ws.Range("L113").Value = 3750
$ws.Range("N113").Value = -8090
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 116437.11
$ws.Range("I132").Value = 89454.5
$ws.Range("J132").Value = 170402.33
$ws.Range("K132").Value = 268363.5
$ws.Range("L132").Value = 511206.99
$ws.Range("M132").Value = -265833.5
$ws.Range("N132").Value = -516266.99
$ws.Range("H134").Value = 4171.727
$ws.Range("I134").Value = 3941.2856
$ws.Range("J134").Value = 4575
$ws.Range("K134").Value = 11823.8568
$ws.Range("L134").Value = 13725
$ws.Range("M134").Value = -9288.856800000001
$ws.Range("N134").Value = -18795
$ws.Range("H136").Value = 3119.2856
$ws.Range("I136").Value = 3015.3157
$ws.Range("J136").Value = 4107
$ws.Range("K136").Value = 9045.947100000001
$ws.Range("L136").Value = 12321
$ws.Range("M136").Value = -6495.947100000001
$ws.Range("N136").Value = -17421

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 60.727272
$ws.Range("I2").Value = 78.28570999999999
$ws.Range("J2").Value = 30
$ws.Range("K2").Value = 469.71426
$ws.Range("L2").Value = 180
$ws.Range("M2").Value = -356.71426
$ws.Range("N2").Value = -406
$ws.Range("H75").Value = 39742
$ws.Range("J75").Value = 79000
$ws.Range("L75").Value = 237000
$ws.Range("N75").Value = -238996
$ws.Range("H78").Value = 39742
$ws.Range("J78").Value = 79000
$ws.Range("L78").Value = 711000
$ws.Range("N78").Value = -720984
$ws.Range("H80").Value = 4499.9
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 4874.875
$ws.Range("K80").Value = 9000
$ws.Range("L80").Value = 14624.625
$ws.Range("M80").Value = -8064
$ws.Range("N80").Value = -16496.625
$ws.Range("H83").Value = 4499.9
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 4874.875
$ws.Range("K83").Value = 27000
$ws.Range("L83").Value = 43873.875
$ws.Range("M83").Value = -22320
$ws.Range("N83").Value = -53233.875
$ws.Range("H123").Value = 854.6
$ws.Range("I123").Value = 318.25
$ws.Range("J123").Value = 3000
$ws.Range("K123").Value = 954.75
$ws.Range("L123").Value = 9000
$ws.Range("M123").Value = 1495.25
$ws.Range("N123").Value = -13900
$ws.Range("H127").Value = 2987.1667
$ws.Range("J127").Value = 2987.1667
$ws.Range("L127").Value = 8961.500100000001
$ws.Range("N127").Value = -18881.5001
$ws.Range("H129").Value = 1012.4167
$ws.Range("I129").Value = 559.2222
$ws.Range("J129").Value = 2372
$ws.Range("K129").Value = 1677.6666
$ws.Range("L129").Value = 7116
$ws.Range("M129").Value = 3322.3334
$ws.Range("N129").Value = -17116
$ws.Range("H131").Value = 1691.6364
$ws.Range("I131").Value = 1115
$ws.Range("J131").Value = 1819.7778
$ws.Range("K131").Value = 3345
$ws.Range("L131").Value = 5459.3334
$ws.Range("M131").Value = 1695
$ws.Range("N131").Value = -15539.3334

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H62").Value = 135000
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 135000
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 135000
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -136372
$ws.Range("H63").Value = 113993
$ws.Range("I63").Value = 90001
$ws.Range("J63").Value = 137985
$ws.Range("K63").Value = 90001
$ws.Range("L63").Value = 137985
$ws.Range("M63").Value = -89315
$ws.Range("N63").Value = -139357
$ws.Range("H65").Value = 135000
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 135000
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 405000
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -411864
$ws.Range("H66").Value = 113993
$ws.Range("I66").Value = 90001
$ws.Range("J66").Value = 137985
$ws.Range("K66").Value = 270003
$ws.Range("L66").Value = 413955
$ws.Range("M66").Value = -266571
$ws.Range("N66").Value = -420819
$ws.Range("H70").Value = 16149.903
$ws.Range("J70").Value = 5142.857
$ws.Range("L70").Value = 5142.857
$ws.Range("N70").Value = -5682.857
$ws.Range("H73").Value = 16149.903
$ws.Range("J73").Value = 5142.857
$ws.Range("L73").Value = 5142.857
$ws.Range("N73").Value = -7014.857
$ws.Range("H75").Value = 135000
$ws.Range("J75").Value = 135000
$ws.Range("L75").Value = 135000
$ws.Range("N75").Value = -136748
$ws.Range("H78").Value = 135000
$ws.Range("J78").Value = 135000
$ws.Range("L78").Value = 405000
$ws.Range("N78").Value = -413736
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H82").Value = 118999.336
$ws.Range("J82").Value = 156999
$ws.Range("L82").Value = 156999
$ws.Range("N82").Value = -157765
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H85").Value = 118999.336
$ws.Range("J85").Value = 156999
$ws.Range("L85").Value = 156999
$ws.Range("N85").Value = -159651
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").ClearContents()
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H96").Value = 63328
$ws.Range("J96").Value = 63328
$ws.Range("L96").Value = 63328
$ws.Range("N96").Value = -68820
$ws.Range("H97").Value = 985.7143
$ws.Range("I97").Value = 860
$ws.Range("J97").Value = 3500
$ws.Range("K97").Value = 860
$ws.Range("L97").Value = 3500
$ws.Range("M97").Value = -364
$ws.Range("N97").Value = -4492
$ws.Range("H102").Value = 1404.6818
$ws.Range("I102").Value = 1281.0952
$ws.Range("J102").Value = 4000
$ws.Range("K102").Value = 1281.0952
$ws.Range("L102").Value = 4000
$ws.Range("M102").Value = 340.9048
$ws.Range("N102").Value = -7244
$ws.Range("H122").Value = 6878.25
$ws.Range("I122").Value = 6501.6665
$ws.Range("J122").Value = 8008
$ws.Range("K122").Value = 19504.9995
$ws.Range("L122").Value = 24024
$ws.Range("M122").Value = -17054.9995
$ws.Range("N122").Value = -28924
$ws.Range("H125").Value = 84629.60000000001
$ws.Range("J125").Value = 84629.60000000001
$ws.Range("L125").Value = 84629.60000000001
$ws.Range("N125").Value = -89549.60000000001
$ws.Range("H126").Value = 3500
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 3500
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 10500
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -15440
$ws.Range("H127").Value = 91261
$ws.Range("J127").Value = 91261
$ws.Range("L127").Value = 91261
$ws.Range("N127").Value = -101181
$ws.Range("H132").Value = 7433.875
$ws.Range("I132").Value = 7433.875
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 22301.625
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -19771.625
$ws.Range("N132").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10795.583
$ws.Range("I7").Value = 16424.75
$ws.Range("J7").Value = 7981
$ws.Range("K7").Value = 16424.75
$ws.Range("L7").Value = 7981
$ws.Range("M7").Value = -16312.75
$ws.Range("N7").Value = -8205
$ws.Range("H40").Value = 37047260
$ws.Range("I40").Value = 66673668
$ws.Range("K40").Value = 66673668
$ws.Range("M40").Value = -66673532
$ws.Range("H61").Value = 2899
$ws.Range("I61").Value = 1689.8636
$ws.Range("J61").Value = 6699.143
$ws.Range("K61").Value = 1689.8636
$ws.Range("L61").Value = 6699.143
$ws.Range("M61").Value = -1487.8636
$ws.Range("N61").Value = -7103.143
$ws.Range("H100").Value = 3642.3333
$ws.Range("I100").Value = 2914.875
$ws.Range("J100").Value = 4473.7144
$ws.Range("K100").Value = 2914.875
$ws.Range("L100").Value = 4473.7144
$ws.Range("M100").Value = -2373.875
$ws.Range("N100").Value = -5555.7144
$ws.Range("H111").Value = 79989
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 79989
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 79989
$ws.Range("M111").ClearContents()
$ws.Range("N111").Value = -88169
$ws.Range("H113").Value = 2899
$ws.Range("I113").Value = 1689.8636
$ws.Range("J113").Value = 6699.143
$ws.Range("K113").Value = 1689.8636
$ws.Range("L113").Value = 6699.143
$ws.Range("M113").Value = 480.1364000000001
$ws.Range("N113").Value = -11039.143
$ws.Range("H126").Value = 10795.583
$ws.Range("I126").Value = 16424.75
$ws.Range("J126").Value = 7981
$ws.Range("K126").Value = 49274.25
$ws.Range("L126").Value = 23943
$ws.Range("M126").Value = -46804.25
$ws.Range("N126").Value = -28883
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 487
$ws.Range("I107").Value = 451.66666
$ws.Range("J107").Value = 699
$ws.Range("K107").Value = 1354.99998
$ws.Range("L107").Value = 2097
$ws.Range("M107").Value = 565.0000199999999
$ws.Range("N107").Value = -5937
$ws.Range("H122").Value = 125005790
$ws.Range("I122").Value = 333335040
$ws.Range("J122").Value = 8244.200000000001
$ws.Range("K122").Value = 1000005120
$ws.Range("L122").Value = 24732.6
$ws.Range("M122").Value = -1000002670
$ws.Range("N122").Value = -29632.6
$ws.Range("H126").Value = 9634.333000000001
$ws.Range("I126").Value = 10361.2
$ws.Range("J126").Value = 6000
$ws.Range("K126").Value = 31083.6
$ws.Range("L126").Value = 18000
$ws.Range("M126").Value = -28613.6
$ws.Range("N126").Value = -22940
$ws.Range("H127").Value = 82979
$ws.Range("J127").Value = 82979
$ws.Range("L127").Value = 82979
$ws.Range("N127").Value = -92899
$ws.Range("H132").Value = 4394.2964
$ws.Range("I132").Value = 3669.6155
$ws.Range("J132").Value = 5067.2144
$ws.Range("K132").Value = 11008.8465
$ws.Range("L132").Value = 15201.6432
$ws.Range("M132").Value = -8478.8465
$ws.Range("N132").Value = -20261.6432
$ws.Range("H136").Value = 30886.885
$ws.Range("I136").Value = 1581.238
$ws.Range("K136").Value = 4743.714
$ws.Range("M136").Value = -2193.714
